$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "261.30"
Set-TextValue $ws "E2" "0.37%"
Set-TextValue $ws "D3" "26.57"
Set-TextValue $ws "E3" "-2.64%"
Set-TextValue $ws "D4" "4.705"
Set-TextValue $ws "E4" "0.23%"
Set-TextValue $ws "D5" "0.06081"
Set-TextValue $ws "E5" "-0.44%"
Set-TextValue $ws "D6" "6.706"
Set-TextValue $ws "E6" "0.73%"
Set-TextValue $ws "D7" "0.8525"
Set-TextValue $ws "E7" "0.29%"
Set-TextValue $ws "D8" "0.9125"
Set-TextValue $ws "E8" "-0.95%"
Set-TextValue $ws "D9" "0.1402"
Set-TextValue $ws "E9" "0.27%"
Set-TextValue $ws "D10" "0.05133"
Set-TextValue $ws "E10" "7.82%"
Set-TextValue $ws "D11" "0.07091"
Set-TextValue $ws "E11" "0.07%"
Set-TextValue $ws "D12" "0.03115"
Set-TextValue $ws "E12" "1.26%"
Set-TextValue $ws "D13" "0.09044"
Set-TextValue $ws "E13" "-0.10%"
Set-TextValue $ws "D14" "0.001549"
Set-TextValue $ws "E14" "1.05%"
Set-TextValue $ws "D15" "0.0006180"
Set-TextValue $ws "E15" "1.29%"
Set-TextValue $ws "D16" "0.006144"
Set-TextValue $ws "E16" "-0.03%"
Set-TextValue $ws "D17" "3.451"
Set-TextValue $ws "E17" "0.00%"
Set-TextValue $ws "D18" "3.172"
Set-TextValue $ws "E18" "0.73%"
Set-TextValue $ws "D19" "2.167"
Set-TextValue $ws "E19" "0.21%"
Set-TextValue $ws "D21" "0.1301"
Set-TextValue $ws "E21" "-0.25%"
Set-TextValue $ws "D22" "4.103"
Set-TextValue $ws "E22" "0.17%"
Set-TextValue $ws "D23" "0.04239"
Set-TextValue $ws "E23" "0.10%"
Set-TextValue $ws "D24" "0.001179"
Set-TextValue $ws "E24" "-3.57%"
Set-TextValue $ws "D25" "0.004044"
Set-TextValue $ws "E25" "6.28%"
Set-TextValue $ws "E26" "0.07%"
Set-TextValue $ws "E27" "4.12%"
Set-TextValue $ws "D40" "0.03979"
Set-TextValue $ws "E40" "3.14%"
Set-TextValue $ws "D41" "0.1111"
Set-TextValue $ws "E41" "-0.06%"
Set-TextValue $ws "D42" "0.004183"
Set-TextValue $ws "E42" "2.56%"
Set-TextValue $ws "B43" "LocalTraders"
Set-TextValue $ws "C43" "https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct"
Set-TextValue $ws "D43" "0.01390"
Set-TextValue $ws "E43" "-14.73%"
Set-TextValue $ws "B44" "CEJI"
Set-TextValue $ws "C44" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D44" "0.002071"
Set-TextValue $ws "E44" "-6.58%"
Set-TextValue $ws "D45" "0.00005117"
Set-TextValue $ws "E45" "-0.80%"
Set-TextValue $ws "E46" "0.07%"
Set-TextValue $ws "D48" "0.2577"
Set-TextValue $ws "E48" "90.17%"
Set-TextValue $ws "E49" "0.07%"
Set-TextValue $ws "E50" "0.07%"
